$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Strings" -> "Strings " + "(strings suck never forget)" as two runs
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Text = "Strings "

# Create the second bit of text in its own paragraph (so it lands in a
# fresh run), then delete the paragraph mark joining them so both runs end
# up together inside paragraph 2, but remain distinct <w:r> elements.
$r2.InsertParagraphAfter() | Out-Null
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "(strings suck never forget)"

$p2 = $d.Paragraphs(2)
$joinMark = $d.Range($p2.Range.End - 1, $p2.Range.End)
$joinMark.Delete()

# ---------------------------------------------------------------------------
# 2) Append the new "Hashing" section after the last paragraph
#    ("El kmp te permite ...")
# ---------------------------------------------------------------------------
function Add-Para([string]$text) {
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter() | Out-Null
    $new = $d.Paragraphs($d.Paragraphs.Count)
    if ($text -ne $null -and $text -ne "") {
        $new.Range.Text = $text
    }
    return $new
}

Add-Para ""                                                                  | Out-Null
Add-Para "Hashing:"                                                          | Out-Null
Add-Para "(puede fallar!!)"                                                  | Out-Null
Add-Para "Es muy rapido"                                                     | Out-Null
Add-Para ""                                                                  | Out-Null
Add-Para "Ver algoritmo de Rabin-Karp"                                       | Out-Null
Add-Para "Ejemplo de aplicacion: encontrar el carácter donde dos strings difieren" | Out-Null

# That paragraph needs a trailing " " in its own separate run, same trick
# as with "Strings " above.
$pEj = $d.Paragraphs($d.Paragraphs.Count)
$pEj.Range.InsertParagraphAfter() | Out-Null
$pSpace = $d.Paragraphs($d.Paragraphs.Count)
$pSpace.Range.Text = " "
$pEj = $d.Paragraphs($d.Paragraphs.Count - 1)
$joinMark2 = $d.Range($pEj.Range.End - 1, $pEj.Range.End)
$joinMark2.Delete()

Add-Para ""                                                                  | Out-Null
Add-Para "Suffix array: O(n logn) "                                          | Out-Null
Add-Para "Suffix tree: trie de suffix arrays"                                | Out-Null
